{"js": "async (context) => {\n  // Positional mapping of the old two-digit multiplication expressions\n  // to their new replacements, in document order.\n  const replacements = [\n    [\"78\u00d716=\", \"18\u00d777=\"],\n    [\"31\u00d758=\", \"31\u00d779=\"],\n    [\"13\u00d753=\", \"94\u00d750=\"],\n    [\"60\u00d759=\", \"21\u00d722=\"],\n    [\"81\u00d775=\", \"96\u00d798=\"],\n    [\"24\u00d713=\", \"59\u00d718=\"],\n    [\"60\u00d749=\", \"25\u00d729=\"],\n    [\"73\u00d791=\", \"93\u00d739=\"],\n    [\"22\u00d751=\", \"40\u00d746=\"],\n    [\"82\u00d755=\", \"59\u00d767=\"],\n    [\"75\u00d739=\", \"83\u00d799=\"],\n    [\"36\u00d764=\", \"82\u00d752=\"],\n    [\"15\u00d774=\", \"20\u00d713=\"],\n    [\"56\u00d789=\", \"67\u00d719=\"],\n    [\"85\u00d728=\", \"38\u00d796=\"],\n    [\"74\u00d718=\", \"79\u00d788=\"],\n    [\"97\u00d720=\", \"23\u00d762=\"],\n    [\"13\u00d770=\", \"94\u00d770=\"],\n    [\"71\u00d768=\", \"58\u00d723=\"],\n    [\"43\u00d783=\", \"62\u00d712=\"],\n    [\"88\u00d765=\", \"53\u00d749=\"],\n    [\"52\u00d712=\", \"26\u00d758=\"],\n    [\"30\u00d757=\", \"38\u00d745=\"],\n    [\"84\u00d796=\", \"70\u00d753=\"],\n    [\"96\u00d765=\", \"13\u00d712=\"],\n  ];\n\n  const body = context.document.body;\n\n  for (const [oldText, newText] of replacements) {\n    const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n    results.load(\"items\");\n    await context.sync();\n\n    if (results.items.length === 0) {\n      throw new Error(`Could not find text to replace: ${oldText}`);\n    }\n\n    // Replace only the first occurrence to preserve a strict 1:1 positional\n    // mapping (the source values are all unique within the document).\n    results.items[0].insertText(newText, \"Replace\");\n  }\n\n  await context.sync();\n};\n", "ps1": "$d = $word.ActiveDocument\n\n# Positional mapping of the old two-digit multiplication expressions\n# to their new replacements, in document order. Each old value is\n# unique within the document, so a single wdReplaceAll (max count 1)\n# find/replace per pair is unambiguous.\n$replacements = @(\n    @{ Old = \"78\u00d716=\"; New = \"18\u00d777=\" },\n    @{ Old = \"31\u00d758=\"; New = \"31\u00d779=\" },\n    @{ Old = \"13\u00d753=\"; New = \"94\u00d750=\" },\n    @{ Old = \"60\u00d759=\"; New = \"21\u00d722=\" },\n    @{ Old = \"81\u00d775=\"; New = \"96\u00d798=\" },\n    @{ Old = \"24\u00d713=\"; New = \"59\u00d718=\" },\n    @{ Old = \"60\u00d749=\"; New = \"25\u00d729=\" },\n    @{ Old = \"73\u00d791=\"; New = \"93\u00d739=\" },\n    @{ Old = \"22\u00d751=\"; New = \"40\u00d746=\" },\n    @{ Old = \"82\u00d755=\"; New = \"59\u00d767=\" },\n    @{ Old = \"75\u00d739=\"; New = \"83\u00d799=\" },\n    @{ Old = \"36\u00d764=\"; New = \"82\u00d752=\" },\n    @{ Old = \"15\u00d774=\"; New = \"20\u00d713=\" },\n    @{ Old = \"56\u00d789=\"; New = \"67\u00d719=\" },\n    @{ Old = \"85\u00d728=\"; New = \"38\u00d796=\" },\n    @{ Old = \"74\u00d718=\"; New = \"79\u00d788=\" },\n    @{ Old = \"97\u00d720=\"; New = \"23\u00d762=\" },\n    @{ Old = \"13\u00d770=\"; New = \"94\u00d770=\" },\n    @{ Old = \"71\u00d768=\"; New = \"58\u00d723=\" },\n    @{ Old = \"43\u00d783=\"; New = \"62\u00d712=\" },\n    @{ Old = \"88\u00d765=\"; New = \"53\u00d749=\" },\n    @{ Old = \"52\u00d712=\"; New = \"26\u00d758=\" },\n    @{ Old = \"30\u00d757=\"; New = \"38\u00d745=\" },\n    @{ Old = \"84\u00d796=\"; New = \"70\u00d753=\" },\n    @{ Old = \"96\u00d765=\"; New = \"13\u00d712=\" }\n)\n\nforeach ($pair in $replacements) {\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $found = $find.Execute(\n        $pair.Old,   # FindText\n        $true,       # MatchCase\n        $true,       # MatchWholeWord\n        $false,      # MatchWildcards\n        $false,      # MatchSoundsLike\n        $false,      # MatchAllWordForms\n        $true,       # Forward\n        1,           # Wrap (wdFindContinue)\n        $false,      # Format\n        $pair.New,   # ReplaceWith\n        1            # Replace (wdReplaceOne)\n    )\n    if (-not $found) {\n        throw \"Could not find text to replace: $($pair.Old)\"\n    }\n}\n"}
